# Remove the "is_locked_lbl" and "is_enabled_lbl" columns (D and E) from the
# detail-view header row. Deleting column D twice removes both, since after
# the first delete the former "is_enabled_lbl" column (E) shifts into D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Delete()
$ws.Columns("D").Delete()

# Append a new "tenant_id_lbl" column at the end (now column F), replacing
# the removed columns' spot in the header row.
$ws.Cells.Item(1, 6).Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'
